$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 103.8666666666667
$ws.Range("E2").Value = 28.9
$ws.Range("F2").ClearContents()

# Row 3
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# Row 4
$ws.Range("C4").Value = 103.55
$ws.Range("E4").Value = 28.25
$ws.Range("F4").ClearContents()

# Row 5
$ws.Range("F5").ClearContents()

# Row 6
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

# Row 7
$ws.Range("C7").Value = 103.4
$ws.Range("E7").Value = 29.2
$ws.Range("F7").ClearContents()

# Row 8
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
